# Update the ObjTables metadata (date + objTablesVersion) stamped into the
# title cell (A1, or A1/A2 for the table-of-contents sheet) of each sheet,
# and fill in the newly-added "Verbose name" column (E) for the Schema
# sheet's attribute rows. All three worksheets are protected with only
# unlocked cells editable by default, so the title cells (locked) require a
# temporary Unprotect/Protect round trip while the Schema sheet's E4:E7
# cells (unlocked) can be written directly.

$wb = $excel.ActiveWorkbook

# --- Sheet "!!_Table of contents" ---
$wsToc = $wb.Worksheets.Item("!!_Table of contents")
$wsToc.Unprotect()
$wsToc.Range("A1").Value = "!!!ObjTables objTablesVersion='1.0.0' date='2020-03-11 23:55:18'"
$wsToc.Range("A2").Value = "!!ObjTables type='TableOfContents' tableFormat='row' description='Table of contents' date='2020-03-11 23:55:18' objTablesVersion='1.0.0'"
$wsToc.Protect()

# --- Sheet "!!_Schema" ---
$wsSchema = $wb.Worksheets.Item("!!_Schema")
$wsSchema.Unprotect()
$wsSchema.Range("A1").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='2020-03-11 23:55:18' objTablesVersion='1.0.0'"
$wsSchema.Protect()

# E4:E7 ("Verbose name" column for the amount/category/date/payee attribute
# rows) use the unlocked data-entry style, so they can be set while the
# sheet stays protected.
$wsSchema.Range("E4").Value = "Amount"
$wsSchema.Range("E5").Value = "Category"
$wsSchema.Range("E6").Value = "Date"
$wsSchema.Range("E7").Value = "Payee"

# --- Sheet "!!Transaction" ---
$wsData = $wb.Worksheets.Item("!!Transaction")
$wsData.Unprotect()
$wsData.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' class='Transaction' name='Transaction' description='Stores transactions' date='2020-03-11 23:55:18' objTablesVersion='1.0.0'"
$wsData.Protect()
